# Auto-generated script applying scheduled market-data refresh values
# to the Carbuncle_Profits leve-profit tables across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1392781.6
$ws.Range("J17").Value = 1392781.6
$ws.Range("L17").Value = 4178344.8
$ws.Range("N17").Value = -4178680.8
$ws.Range("H18").Value = 298
$ws.Range("J18").Value = 500
$ws.Range("L18").Value = 500
$ws.Range("N18").Value = -1068
$ws.Range("H43").Value = 1435.3529
$ws.Range("I43").Value = 990.3333
$ws.Range("J43").Value = 1530.7142
$ws.Range("K43").Value = 990.3333
$ws.Range("L43").Value = 1530.7142
$ws.Range("M43").Value = -921.3333
$ws.Range("N43").Value = -1668.7142
$ws.Range("H51").Value = 8491.833000000001
$ws.Range("J51").Value = 8809.272000000001
$ws.Range("L51").Value = 8809.272000000001
$ws.Range("N51").Value = -9777.272000000001
$ws.Range("H86").Value = 36838.65
$ws.Range("I86").Value = 26138.625
$ws.Range("J86").Value = 46349.777
$ws.Range("K86").Value = 26138.625
$ws.Range("L86").Value = 46349.777
$ws.Range("M86").Value = -25015.625
$ws.Range("N86").Value = -48595.777
$ws.Range("H87").Value = 15570.611
$ws.Range("J87").Value = 15570.611
$ws.Range("L87").Value = 15570.611
$ws.Range("N87").Value = -18066.611
$ws.Range("H89").Value = 36838.65
$ws.Range("I89").Value = 26138.625
$ws.Range("J89").Value = 46349.777
$ws.Range("K89").Value = 130693.125
$ws.Range("L89").Value = 231748.885
$ws.Range("M89").Value = -125077.125
$ws.Range("N89").Value = -242980.885
$ws.Range("H90").Value = 15570.611
$ws.Range("J90").Value = 15570.611
$ws.Range("L90").Value = 46711.833
$ws.Range("N90").Value = -59191.833
$ws.Range("H96").Value = 62502956
$ws.Range("I96").Value = 125002910
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 375008730
$ws.Range("L96").Value = 9000
$ws.Range("M96").Value = -375007357
$ws.Range("N96").Value = -11746
$ws.Range("H100").Value = 3383.2693
$ws.Range("I100").Value = 2966.5908
$ws.Range("J100").Value = 5675
$ws.Range("K100").Value = 2966.5908
$ws.Range("L100").Value = 5675
$ws.Range("M100").Value = -2425.5908
$ws.Range("N100").Value = -6757
$ws.Range("H112").Value = 1247.8966
$ws.Range("J112").Value = 1324.5834
$ws.Range("L112").Value = 3973.7502
$ws.Range("N112").Value = -6189.7502
$ws.Range("H138").Value = 4280.879
$ws.Range("J138").Value = 5187.8438
$ws.Range("L138").Value = 15563.5314
$ws.Range("N138").Value = -25843.5314
$ws.Range("H141").Value = 2986.3333
$ws.Range("I141").Value = 1464
$ws.Range("J141").Value = 38000
$ws.Range("K141").Value = 4392
$ws.Range("L141").Value = 114000
$ws.Range("M141").Value = 788
$ws.Range("N141").Value = -124360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 27666.666
$ws.Range("J7").Value = 27666.666
$ws.Range("L7").Value = 27666.666
$ws.Range("N7").Value = -27894.666
$ws.Range("H32").Value = 5182.59
$ws.Range("I32").Value = 5189.485
$ws.Range("K32").Value = 5189.485
$ws.Range("M32").Value = -4902.485
$ws.Range("H93").Value = 31448
$ws.Range("J93").Value = 31448
$ws.Range("L93").Value = 31448
$ws.Range("N93").Value = -36440
$ws.Range("H97").Value = 621.53845
$ws.Range("I97").Value = 315.29413
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 315.29413
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = 180.70587
$ws.Range("N97").Value = -2192
$ws.Range("H102").Value = 1875
$ws.Range("I102").Value = 1875
$ws.Range("K102").Value = 1875
$ws.Range("M102").Value = -253
$ws.Range("H139").Value = 63310.445
$ws.Range("J139").Value = 63310.445
$ws.Range("L139").Value = 63310.445
$ws.Range("N139").Value = -73590.44500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 53639.332
$ws.Range("J70").Value = 53639.332
$ws.Range("L70").Value = 53639.332
$ws.Range("N70").Value = -54225.332
$ws.Range("H73").Value = 53639.332
$ws.Range("J73").Value = 53639.332
$ws.Range("L73").Value = 53639.332
$ws.Range("N73").Value = -55667.332
$ws.Range("H86").Value = 3092.7856
$ws.Range("I86").Value = 3111
$ws.Range("J86").Value = 3060
$ws.Range("K86").Value = 3111
$ws.Range("L86").Value = 3060
$ws.Range("M86").Value = -1988
$ws.Range("N86").Value = -5306
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496
$ws.Range("H89").Value = 3092.7856
$ws.Range("I89").Value = 3111
$ws.Range("J89").Value = 3060
$ws.Range("K89").Value = 15555
$ws.Range("L89").Value = 15300
$ws.Range("M89").Value = -9939
$ws.Range("N89").Value = -26532
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480
$ws.Range("H94").Value = 521.1053000000001
$ws.Range("I94").Value = 505.94116
$ws.Range("K94").Value = 505.94116
$ws.Range("M94").Value = -54.94116000000002
$ws.Range("H99").Value = 1821.9445
$ws.Range("I99").Value = 1356.0714
$ws.Range("J99").Value = 3452.5
$ws.Range("K99").Value = 1356.0714
$ws.Range("L99").Value = 3452.5
$ws.Range("M99").Value = 141.9286
$ws.Range("N99").Value = -6448.5
$ws.Range("H140").Value = 57331.668
$ws.Range("J140").Value = 57331.668
$ws.Range("L140").Value = 57331.668
$ws.Range("N140").Value = -67691.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 4500
$ws.Range("J17").Value = 8000
$ws.Range("L17").Value = 8000
$ws.Range("N17").Value = -8348
$ws.Range("H25").Value = 5800
$ws.Range("I25").Value = 5800
$ws.Range("K25").Value = 5800
$ws.Range("M25").Value = -5626
$ws.Range("H31").Value = 3670.347
$ws.Range("I31").Value = 1881.5916
$ws.Range("J31").Value = 8374.111000000001
$ws.Range("K31").Value = 1881.5916
$ws.Range("L31").Value = 8374.111000000001
$ws.Range("M31").Value = -1586.5916
$ws.Range("N31").Value = -8964.111000000001
$ws.Range("H34").Value = 3670.347
$ws.Range("I34").Value = 1881.5916
$ws.Range("J34").Value = 8374.111000000001
$ws.Range("K34").Value = 1881.5916
$ws.Range("L34").Value = 8374.111000000001
$ws.Range("M34").Value = -1679.5916
$ws.Range("N34").Value = -8778.111000000001
$ws.Range("H41").Value = 14619.667
$ws.Range("I41").Value = 7029.5
$ws.Range("J41").Value = 29800
$ws.Range("K41").Value = 7029.5
$ws.Range("L41").Value = 29800
$ws.Range("M41").Value = -6601.5
$ws.Range("N41").Value = -30656

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 780.41
$ws.Range("I131").Value = 361.53845
$ws.Range("K131").Value = 1084.61535
$ws.Range("M131").Value = 3955.38465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2152.2273
$ws.Range("I80").Value = 1924.75
$ws.Range("J80").Value = 2202.7778
$ws.Range("K80").Value = 1924.75
$ws.Range("L80").Value = 2202.7778
$ws.Range("M80").Value = -926.75
$ws.Range("N80").Value = -4198.7778
$ws.Range("H83").Value = 2152.2273
$ws.Range("I83").Value = 1924.75
$ws.Range("J83").Value = 2202.7778
$ws.Range("K83").Value = 9623.75
$ws.Range("L83").Value = 11013.889
$ws.Range("M83").Value = -4631.75
$ws.Range("N83").Value = -20997.889
$ws.Range("H97").Value = 1935.3846
$ws.Range("J97").Value = 1300
$ws.Range("L97").Value = 1300
$ws.Range("N97").Value = -2292
$ws.Range("H99").Value = 10374.75
$ws.Range("I99").Value = 10374.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10374.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -8128.75
$ws.Range("N99").ClearContents()
$ws.Range("H138").Value = 50341.7
$ws.Range("J138").Value = 50341.7
$ws.Range("L138").Value = 50341.7
$ws.Range("N138").Value = -60621.7
$ws.Range("H140").Value = 39303.332
$ws.Range("J140").Value = 39303.332
$ws.Range("L140").Value = 39303.332
$ws.Range("N140").Value = -49663.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2431.375
$ws.Range("I82").Value = 2325.2222
$ws.Range("K82").Value = 2325.2222
$ws.Range("M82").Value = -1964.2222
$ws.Range("H85").Value = 2431.375
$ws.Range("I85").Value = 2325.2222
$ws.Range("K85").Value = 2325.2222
$ws.Range("M85").Value = -1077.2222
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 1565
$ws.Range("J93").Value = 1966.6666
$ws.Range("L93").Value = 1966.6666
$ws.Range("N93").Value = -4462.6666
$ws.Range("H99").Value = 34242.5
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H100").Value = 2211.5557
$ws.Range("I100").Value = 1487.5
$ws.Range("J100").Value = 8004
$ws.Range("K100").Value = 1487.5
$ws.Range("L100").Value = 8004
$ws.Range("M100").Value = -946.5
$ws.Range("N100").Value = -9086

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 42666.668
$ws.Range("J51").Value = 42666.668
$ws.Range("L51").Value = 42666.668
$ws.Range("N51").Value = -43686.668
$ws.Range("H96").Value = 1779.4
$ws.Range("I96").Value = 1836.75
$ws.Range("J96").Value = 1550
$ws.Range("K96").Value = 1836.75
$ws.Range("L96").Value = 1550
$ws.Range("M96").Value = -463.75
$ws.Range("N96").Value = -4296

